$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# New lattice multiplication content for each of the 5 rows x 3 columns cells.
# Each cell has 5 lines joined with a Word line-break (Chr(11)):
#   "NN x NN", "  D    D", "  ----", "D|    |", "D|    |"
$cellData = @(
    @(1, 1, "46 x 30", "  3    0", "4|    |", "6|    |"),
    @(1, 2, "84 x 42", "  4    2", "8|    |", "4|    |"),
    @(1, 3, "51 x 26", "  2    6", "5|    |", "1|    |"),
    @(2, 1, "29 x 83", "  8    3", "2|    |", "9|    |"),
    @(2, 2, "75 x 73", "  7    3", "7|    |", "5|    |"),
    @(2, 3, "50 x 91", "  9    1", "5|    |", "0|    |"),
    @(3, 1, "55 x 52", "  5    2", "5|    |", "5|    |"),
    @(3, 2, "94 x 85", "  8    5", "9|    |", "4|    |"),
    @(3, 3, "13 x 51", "  5    1", "1|    |", "3|    |"),
    @(4, 1, "10 x 26", "  2    6", "1|    |", "0|    |"),
    @(4, 2, "23 x 43", "  4    3", "2|    |", "3|    |"),
    @(4, 3, "61 x 67", "  6    7", "6|    |", "1|    |"),
    @(5, 1, "81 x 45", "  4    5", "8|    |", "1|    |"),
    @(5, 2, "75 x 20", "  2    0", "7|    |", "5|    |"),
    @(5, 3, "19 x 15", "  1    5", "1|    |", "9|    |")
)

$lineBreak = [char]11

foreach ($entry in $cellData) {
    $row = $entry[0]
    $col = $entry[1]
    $header = $entry[2]
    $middle = $entry[3]
    $line4 = $entry[4]
    $line5 = $entry[5]

    $newText = $header + $lineBreak + $middle + $lineBreak + "  ----" + $lineBreak + $line4 + $lineBreak + $line5

    $cell = $tbl.Cell($row, $col)
    $cell.Range.Text = $newText
}
